$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555197434459"
$wb.Worksheets.Item(2).Name = "NB_TO-1651255520556097"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555205620975"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555206200964"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555206970956"

# Sheet 1 (GNG) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555197144458.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555197274494.csv"
$ws1.Range("B4").Value = "go_stims-16512555197294445.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255519742461.csv"

# Sheet 2 (NB) updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512555201720972.csv"
$ws2.Range("B3").Value = "OB-1651255520074097.csv"
$ws2.Range("B4").Value = "TB-16512555205300965.csv"
$ws2.Range("B5").Value = "ZB-match_7-16512555200200982.csv"
$ws2.Range("B6").Value = "TB-1651255520231097.csv"
$ws2.Range("B7").Value = "TB-16512555204320958.csv"
$ws2.Range("B8").Value = "ZB-match_5-16512555197684462.csv"
$ws2.Range("B9").Value = "OB-16512555201520972.csv"
$ws2.Range("B10").Value = "ZB-match_6-16512555198464441.csv"

# Sheet 4 (TOL) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555205871003.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255520564098.csv"
$ws4.Range("B4").Value = "MM_stims-16512555206030986.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555205880995.csv"
$ws4.Range("B6").Value = "MM_stims-16512555206190972.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555206041017.csv"

# Sheet 5 (vSAT) updates
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555206250978.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555206501012.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651255520682099.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555206660984.csv"
